$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Range("F1").Value = "total_time"
$ws.Range("H1").Value = "complete"
$ws.Range("G1").Value = "waiting"

# Row 2
$ws.Range("F2").Value = 0.002057

# Row 3
$ws.Range("F3").Value = 0.006406

# Row 4
$ws.Range("E4").Value = 50
$ws.Range("F4").Value = 0.009025
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0

# Row 5
$ws.Range("E5").Value = 60
$ws.Range("F5").Value = 0.016001
$ws.Range("G5").Value = 1333
$ws.Range("H5").Value = 2000

# Row 6
$ws.Range("E6").Value = 60
$ws.Range("F6").Value = 0.023902
$ws.Range("G6").Value = 2500
$ws.Range("H6").Value = 3500

# Row 7
$ws.Range("E7").Value = 70
$ws.Range("F7").Value = 0.026304
$ws.Range("G7").Value = 7666
$ws.Range("H7").Value = 8333

# Row 8
$ws.Range("E8").Value = 70
$ws.Range("F8").Value = 0.054008
$ws.Range("G8").Value = 9400
$ws.Range("H8").Value = 10000

$ws.Range("F11").Select()
